$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Update shortname value on input sheet (B2) from numeric 2425 to text "2425d"
$wsInput.Range("B2").Value = "2425d"

# Update productname value on input sheet (B1) with the "-1st" suffixed name
$wsInput.Range("B1").Value = "2425-RBI-EI-DB-DL-REC-NOCOM-RNI-CTPD-DL-MD-TR-1-DATE-VAR-INST-1st"

# Mirror the productname update onto the output sheet (B1)
$wsOutput.Range("B1").Value = "2425-RBI-EI-DB-DL-REC-NOCOM-RNI-CTPD-DL-MD-TR-1-DATE-VAR-INST-1st"

# Reset selection on input sheet to B1 (no more special top-left / selection at A28)
$wsInput.Range("B1").Select()

# Select B1 on output sheet as well, then make the output sheet the active sheet/tab
$wsOutput.Range("B1").Select()
$wsOutput.Activate()
